# Tabelle strategia 3 complete
# Fill in the measured execution-time values (column C) for the P=2, P=4 and
# P=8 strategy tables. The D (speedup) and E (efficiency) columns already
# hold formulas referencing these cells, so they will recalculate
# automatically once the inputs are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tabella P = 2 (rows 11-15) ---
$ws.Range("C11").Value = 0.83699999999999997
$ws.Range("C12").Value = 3.1480000000000001
$ws.Range("C13").Value = 11.766
$ws.Range("C14").Value = 47.704000000000001
$ws.Range("C15").Value = 198.97559999999999

# --- Tabella P = 4 (rows 19-23) ---
$ws.Range("C19").Value = 0.65459999999999996
$ws.Range("C20").Value = 1.8367
$ws.Range("C21").Value = 6.6443000000000003
$ws.Range("C22").Value = 27.895700000000001
$ws.Range("C23").Value = 138.9153

# --- Tabella P = 8 (rows 28-32) ---
$ws.Range("C28").Value = 0.62470000000000003
$ws.Range("C29").Value = 1.9390000000000001
$ws.Range("C30").Value = 5.5472999999999999
$ws.Range("C31").Value = 21.135000000000002
$ws.Range("C32").Value = 97.0167

# C31 picks up a dedicated number format (#,##0.000) distinct from the
# other cells in that column (which keep the default General format).
$ws.Range("C31").NumberFormat = "#,##0.000"

# Update the saved view/selection state to match the author's session.
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("G44").Select()
